$wb = $excel.ActiveWorkbook

# Add the new sheet at the end of the workbook (after the last existing sheet)
$lastSheet = $wb.Worksheets.Item($wb.Worksheets.Count)
$ws = $wb.Worksheets.Add($null, $lastSheet)
$ws.Name = "所有频道URL"

# All channel URLs, in row order (row 110 and 111 share the same value, as in the source data)
$values = @(
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226016/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226514/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226559/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225588/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226021/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226428/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226007/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226019/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225603/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226010/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225733/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226008/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225734/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225730/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226473/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226565/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225597/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225731/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226568/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226011/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226537/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226571/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225732/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226591/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226476/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225601/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226100/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225765/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225604/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226517/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226543/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226520/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226431/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225602/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225652/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226425/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225783/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225659/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226391/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226498/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225728/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226450/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226251/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226439/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225618/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226409/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226463/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225735/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226345/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226336/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225620/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226341/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226377/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225633/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226522/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226248/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226374/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226549/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226380/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226474/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225626/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226383/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226465/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225639/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226445/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226406/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225623/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226448/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226480/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226320/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225627/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226477/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226501/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225610/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226307/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226430/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226397/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226386/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225613/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226310/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226495/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226344/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226504/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225660/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225619/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226546/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226254/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225634/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226389/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225632/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226454/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225628/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226323/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225655/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226456/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226507/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225624/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226392/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225625/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226457/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225739/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226313/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226510/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226338/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226513/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225740/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226459/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226395/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225638/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226433/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226433/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225635/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226460/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226516/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226444/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225612/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226339/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226342/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226327/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226585/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226582/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226532/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226531/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226525/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226467/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226419/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226421/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226298/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226365/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226239/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226424/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226301/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226427/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226330/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226368/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226242/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226492/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226304/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226371/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226434/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226558/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225760/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226552/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225762/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226526/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226437/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226413/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226555/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226493/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226534/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226496/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226484/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226523/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226499/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226529/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226502/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226451/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226416/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226535/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226422/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226410/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226481/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225657/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225729/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225654/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225656/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225653/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226124/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226147/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226469/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226508/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226472/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226398/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225669/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225716/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225717/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225715/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225714/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225685/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225675/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225741/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225663/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225671/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225743/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225662/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225742/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225677/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225666/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225665/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225670/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225672/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221226505/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225674/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225676/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225668/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225664/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225683/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225680/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225719/index.m3u8',
    '/ottrrs.hl.chinamobile.com/PLTV/88888888/224/3221225673/index.m3u8'
)

for ($i = 0; $i -lt $values.Length; $i++) {
    $ws.Cells.Item($i + 1, 1).Value = $values[$i]
}

# Match formatting of the source sheet: vertical-centered cells, wide column A
$ws.Range("A1:A200").VerticalAlignment = -4108
$ws.Columns.Item(1).ColumnWidth = 73

# Page setup matching the other worksheets (portrait, paper size 9 / A4)
$ws.PageSetup.PaperSize = 9
$ws.PageSetup.Orientation = 1

# Restore the originally active sheet/tab
$wb.Worksheets.Item(1).Activate()
